# Decrement the numeric suffix of each label in column A (rows 2-97):
# q1 -> q0, q2 -> q1, ..., q96 -> q95
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 97; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $current = $cell.Value()
    if ($current -match '^q(\d+)$') {
        $num = [int]$matches[1]
        $cell.Value = "q$($num - 1)"
    }
}
